# Jogos_do_Dia_Betfair_Back_Lay_2025-12-16.xlsx update
# - tweak several odds on rows 2 and 3
# - row 4 now describes a different match (English National League,
#   Truro City vs Wealdstone) replacing the old Swiss Super League fixture
# - rows 5 and 6 (the other two Swiss Super League / English National League
#   rows that used to trail the sheet) are removed entirely

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Row 2 (Thai League 1, Rayong FC vs Ratchaburi) ----
$ws.Range("I2").Value = 2.98
$ws.Range("R2").Value = 1.33
$ws.Range("S2").Value = 2.32
$ws.Range("W2").Value = 1.46

# ---- Row 3 (Portuguese Segunda Liga, Maritimo vs Benfica B) ----
$ws.Range("G3").Value = 2.02
$ws.Range("H3").Value = 4.2
$ws.Range("J3").Value = 3.45
$ws.Range("K3").Value = 3.95
$ws.Range("M3").Value = 1.07
$ws.Range("N3").Value = 3.05
$ws.Range("O3").Value = 1.33
$ws.Range("P3").Value = 1.83
$ws.Range("R3").Value = 1.32
$ws.Range("S3").Value = 3.1
$ws.Range("T3").Value = 1.84
$ws.Range("U3").Value = 1.97
$ws.Range("W3").Value = 1.98
$ws.Range("AB3").Value = 10.5
$ws.Range("AC3").Value = 8.800000000000001
$ws.Range("AG3").Value = 12.5

# ---- Row 4: replaced with English National League, Truro City vs Wealdstone ----
$ws.Range("A4").Value = "English National League"
$ws.Range("C4").Value = "16:45:00"
$ws.Range("D4").Value = "Truro City"
$ws.Range("E4").Value = "Wealdstone"
$ws.Range("F4").Value = 2.7
$ws.Range("G4").Value = 3.1
$ws.Range("H4").Value = 2.48
$ws.Range("I4").Value = 2.8
$ws.Range("J4").Value = 3.35
$ws.Range("K4").Value = 3.95
$ws.Range("N4").Value = 3.85
$ws.Range("O4").Value = 1.3
$ws.Range("P4").Value = 1.97
$ws.Range("Q4").Value = 1.82

# ---- Remove old rows 5 and 6 (Winterthur vs Thun, and the duplicate
#      Truro City vs Wealdstone row that used to live at row 6) ----
$ws.Rows("5:6").Delete()
